$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '30.716.02'
$ws.Range('E2').Value = '  +2.79%  '

$ws.Range('D3').Value = '1.685.54'
$ws.Range('E3').Value = '  +3.04%  '

$ws.Range('E4').Value = '  -0.37%  '

$ws.Range('E5').Value = '  +2.50%  '

$ws.Range('E6').Value = '  +2.07%  '

$ws.Range('E7').Value = '  -0.29%  '

Set-TextValue 'D8' '30.44'
$ws.Range('E8').Value = '  +5.28%  '

$ws.Range('E9').Value = '  +2.52%  '

$ws.Range('E10').Value = '  +3.28%  '

$ws.Range('E11').Value = '  -0.88%  '

$ws.Range('D12').Value = '1.928.34'
$ws.Range('E12').Value = '  +3.09%  '

Set-TextValue 'D13' '10.44'
$ws.Range('E13').Value = '  +12.28%  '

$ws.Range('E14').Value = '  +9.76%  '

$ws.Range('D15').Value = '1.687.44'
$ws.Range('E15').Value = '  +3.11%  '

$ws.Range('E16').Value = '  +2.72%  '

$ws.Range('D17').Value = '30.721.50'
$ws.Range('E17').Value = '  +2.68%  '

Set-TextValue 'D18' '66.42'
$ws.Range('E18').Value = '  +3.34%  '

Set-TextValue 'D19' '245.86'
$ws.Range('E19').Value = '  +0.92%  '

$ws.Range('D20').Value = '0.0₃0715'
$ws.Range('E20').Value = '  +1.95%  '

$ws.Range('E21').Value = '  -0.36%  '

$ws.Range('E22').Value = '  +3.12%  '

$ws.Range('E23').Value = '  +3.60%  '

Set-TextValue 'D24' '2.17'
$ws.Range('E24').Value = '  +1.57%  '

Set-TextValue 'D25' '158.13'
$ws.Range('E25').Value = '  +0.02%  '

Set-TextValue 'D26' '15.89'
$ws.Range('E26').Value = '  +2.07%  '

$ws.Range('E27').Value = '  +0.90%  '

$ws.Range('E28').Value = '  +1.37%  '

Set-TextValue 'D29' '1.00'
$ws.Range('E29').Value = '  -0.47%  '

$ws.Range('E30').Value = '  +2.13%  '

$ws.Range('E31').Value = '  +1.38%  '

$ws.Range('E32').Value = '  +3.03%  '

$ws.Range('D33').Value = '1.513.91'
$ws.Range('E33').Value = '  +6.08%  '

$ws.Range('E34').Value = '  +3.57%  '

Set-TextValue 'D35' '1.75'
$ws.Range('E35').Value = '  +5.53%  '

Set-TextValue 'D36' '84.71'
$ws.Range('E36').Value = '  +10.07%  '

$ws.Range('E37').Value = '  -0.98%  '

$ws.Range('E38').Value = '  +4.86%  '

$ws.Range('E39').Value = '  +6.33%  '

Set-TextValue 'D40' '2.72'
$ws.Range('E40').Value = '  -4.43%  '

$ws.Range('E41').Value = '  -0.15%  '

Set-TextValue 'D42' '0.839'
$ws.Range('E42').Value = '  +1.05%  '

$ws.Range('E43').Value = '  -0.42%  '

$ws.Range('E44').Value = '  +1.36%  '

$ws.Range('E45').Value = '  -1.11%  '

$ws.Range('E46').Value = '  -0.38%  '

Set-TextValue 'D47' '51.92'
$ws.Range('E47').Value = '  -3.27%  '

$ws.Range('D48').Value = '1.820.13'
$ws.Range('E48').Value = '  +2.33%  '

Set-TextValue 'D49' '5.46'
$ws.Range('E49').Value = '  +1.94%  '

Set-TextValue 'D50' '94.70'
$ws.Range('E50').Value = '  +5.78%  '

$ws.Range('E51').Value = '  +2.47%  '
